# Update the LR-pair data table with the new TPM-derived values.
# New sending/target cluster set now also includes "ECs", expanding the
# table from 8 data rows (2 senders x 4 targets) to 12 (3 senders x 4 targets),
# and all numeric metrics are recomputed against the new TPM data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Pthlh"
$ws.Cells.Item(2, 3).Value = "Pth1r"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.2456256666666667
$ws.Cells.Item(2, 8).Value = 0.736877
$ws.Cells.Item(2, 9).Value = 0.03469041475194144
$ws.Cells.Item(2, 10).Value = 0.03469041475194144
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.8566003333333333
$ws.Cells.Item(2, 14).Value = 2.569801
$ws.Cells.Item(2, 15).Value = 0.1153349512295097
$ws.Cells.Item(2, 16).Value = 0.1153349512295097
$ws.Cells.Item(2, 17).Value = 0.2104030279418889
$ws.Cells.Item(2, 18).Value = 1.893627251477
$ws.Cells.Item(2, 19).Value = 0.004001017293546629
$ws.Cells.Item(2, 20).Value = 0.004001017293546629

# Row 3: ECs -> FAPs (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Pthlh"
$ws.Cells.Item(3, 3).Value = "Pth1r"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.2456256666666667
$ws.Cells.Item(3, 8).Value = 0.736877
$ws.Cells.Item(3, 9).Value = 0.03469041475194144
$ws.Cells.Item(3, 10).Value = 0.03469041475194144
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 4.697042333333333
$ws.Cells.Item(3, 14).Value = 14.091127
$ws.Cells.Item(3, 15).Value = 0.6324222946888989
$ws.Cells.Item(3, 16).Value = 0.632422294688899
$ws.Cells.Item(3, 17).Value = 1.153714154486555
$ws.Cells.Item(3, 18).Value = 10.383427390379
$ws.Cells.Item(3, 19).Value = 0.02193899170113243
$ws.Cells.Item(3, 20).Value = 0.02193899170113243

# Row 4: ECs -> MuSCs (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Pthlh"
$ws.Cells.Item(4, 3).Value = "Pth1r"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.2456256666666667
$ws.Cells.Item(4, 8).Value = 0.736877
$ws.Cells.Item(4, 9).Value = 0.03469041475194144
$ws.Cells.Item(4, 10).Value = 0.03469041475194144
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.559594666666667
$ws.Cells.Item(4, 14).Value = 4.678784
$ws.Cells.Item(4, 15).Value = 0.2099879813469643
$ws.Cells.Item(4, 16).Value = 0.2099879813469644
$ws.Cells.Item(4, 17).Value = 0.3830764797297778
$ws.Cells.Item(4, 18).Value = 3.447688317568
$ws.Cells.Item(4, 19).Value = 0.007284570165849134
$ws.Cells.Item(4, 20).Value = 0.007284570165849135

# Row 5: ECs -> Resolving-Mac (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Pthlh"
$ws.Cells.Item(5, 3).Value = "Pth1r"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.2456256666666667
$ws.Cells.Item(5, 8).Value = 0.736877
$ws.Cells.Item(5, 9).Value = 0.03469041475194144
$ws.Cells.Item(5, 10).Value = 0.03469041475194144
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.313829
$ws.Cells.Item(5, 14).Value = 0.941487
$ws.Cells.Item(5, 15).Value = 0.04225477273462707
$ws.Cells.Item(5, 16).Value = 0.04225477273462708
$ws.Cells.Item(5, 17).Value = 0.07708445734433332
$ws.Cells.Item(5, 18).Value = 0.693760116099
$ws.Cells.Item(5, 19).Value = 0.00146583559141324
$ws.Cells.Item(5, 20).Value = 0.00146583559141324

# Row 6: FAPs -> ECs (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Pthlh"
$ws.Cells.Item(6, 3).Value = "Pth1r"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 6.802404
$ws.Cells.Item(6, 8).Value = 20.407212
$ws.Cells.Item(6, 9).Value = 0.9607229540490425
$ws.Cells.Item(6, 10).Value = 0.9607229540490425
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.8566003333333333
$ws.Cells.Item(6, 14).Value = 2.569801
$ws.Cells.Item(6, 15).Value = 0.1153349512295097
$ws.Cells.Item(6, 16).Value = 0.1153349512295097
$ws.Cells.Item(6, 17).Value = 5.826941533868
$ws.Cells.Item(6, 18).Value = 52.442473804812
$ws.Cells.Item(6, 19).Value = 0.1108049350503168
$ws.Cells.Item(6, 20).Value = 0.1108049350503168

# Row 7: FAPs -> FAPs (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Pthlh"
$ws.Cells.Item(7, 3).Value = "Pth1r"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 6.802404
$ws.Cells.Item(7, 8).Value = 20.407212
$ws.Cells.Item(7, 9).Value = 0.9607229540490425
$ws.Cells.Item(7, 10).Value = 0.9607229540490425
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.697042333333333
$ws.Cells.Item(7, 14).Value = 14.091127
$ws.Cells.Item(7, 15).Value = 0.6324222946888989
$ws.Cells.Item(7, 16).Value = 0.632422294688899
$ws.Cells.Item(7, 17).Value = 31.951179556436
$ws.Cells.Item(7, 18).Value = 287.560616007924
$ws.Cells.Item(7, 19).Value = 0.607582615159993
$ws.Cells.Item(7, 20).Value = 0.6075826151599931

# Row 8: FAPs -> MuSCs (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Pthlh"
$ws.Cells.Item(8, 3).Value = "Pth1r"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 6.802404
$ws.Cells.Item(8, 8).Value = 20.407212
$ws.Cells.Item(8, 9).Value = 0.9607229540490425
$ws.Cells.Item(8, 10).Value = 0.9607229540490425
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.559594666666667
$ws.Cells.Item(8, 14).Value = 4.678784
$ws.Cells.Item(8, 15).Value = 0.2099879813469643
$ws.Cells.Item(8, 16).Value = 0.2099879813469644
$ws.Cells.Item(8, 17).Value = 10.608992998912
$ws.Cells.Item(8, 18).Value = 95.48093699020801
$ws.Cells.Item(8, 19).Value = 0.2017402737544508
$ws.Cells.Item(8, 20).Value = 0.2017402737544508

# Row 9: FAPs -> Resolving-Mac (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Pthlh"
$ws.Cells.Item(9, 3).Value = "Pth1r"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 6.802404
$ws.Cells.Item(9, 8).Value = 20.407212
$ws.Cells.Item(9, 9).Value = 0.9607229540490425
$ws.Cells.Item(9, 10).Value = 0.9607229540490425
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.313829
$ws.Cells.Item(9, 14).Value = 0.941487
$ws.Cells.Item(9, 15).Value = 0.04225477273462707
$ws.Cells.Item(9, 16).Value = 0.04225477273462708
$ws.Cells.Item(9, 17).Value = 2.134791644916
$ws.Cells.Item(9, 18).Value = 19.213124804244
$ws.Cells.Item(9, 19).Value = 0.04059513008428186
$ws.Cells.Item(9, 20).Value = 0.04059513008428187

# Row 10: MuSCs -> ECs (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Pthlh"
$ws.Cells.Item(10, 3).Value = "Pth1r"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.03247566666666667
$ws.Cells.Item(10, 8).Value = 0.097427
$ws.Cells.Item(10, 9).Value = 0.004586631199016116
$ws.Cells.Item(10, 10).Value = 0.004586631199016115
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.8566003333333333
$ws.Cells.Item(10, 14).Value = 2.569801
$ws.Cells.Item(10, 15).Value = 0.1153349512295097
$ws.Cells.Item(10, 16).Value = 0.1153349512295097
$ws.Cells.Item(10, 17).Value = 0.02781866689188889
$ws.Cells.Item(10, 18).Value = 0.250368002027
$ws.Cells.Item(10, 19).Value = 0.0005289988856462712
$ws.Cells.Item(10, 20).Value = 0.0005289988856462712

# Row 11: MuSCs -> FAPs (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Pthlh"
$ws.Cells.Item(11, 3).Value = "Pth1r"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.03247566666666667
$ws.Cells.Item(11, 8).Value = 0.097427
$ws.Cells.Item(11, 9).Value = 0.004586631199016116
$ws.Cells.Item(11, 10).Value = 0.004586631199016115
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 4.697042333333333
$ws.Cells.Item(11, 14).Value = 14.091127
$ws.Cells.Item(11, 15).Value = 0.6324222946888989
$ws.Cells.Item(11, 16).Value = 0.632422294688899
$ws.Cells.Item(11, 17).Value = 0.1525395811365556
$ws.Cells.Item(11, 18).Value = 1.372856230229
$ws.Cells.Item(11, 19).Value = 0.002900687827773468
$ws.Cells.Item(11, 20).Value = 0.002900687827773468

# Row 12: MuSCs -> MuSCs (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Pthlh"
$ws.Cells.Item(12, 3).Value = "Pth1r"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.03247566666666667
$ws.Cells.Item(12, 8).Value = 0.097427
$ws.Cells.Item(12, 9).Value = 0.004586631199016116
$ws.Cells.Item(12, 10).Value = 0.004586631199016115
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 1.559594666666667
$ws.Cells.Item(12, 14).Value = 4.678784
$ws.Cells.Item(12, 15).Value = 0.2099879813469643
$ws.Cells.Item(12, 16).Value = 0.2099879813469644
$ws.Cells.Item(12, 17).Value = 0.05064887652977778
$ws.Cells.Item(12, 18).Value = 0.455839888768
$ws.Cells.Item(12, 19).Value = 0.0009631374266644007
$ws.Cells.Item(12, 20).Value = 0.0009631374266644007

# Row 13: MuSCs -> Resolving-Mac (ligand Pthlh / receptor Pth1r)
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Pthlh"
$ws.Cells.Item(13, 3).Value = "Pth1r"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.03247566666666667
$ws.Cells.Item(13, 8).Value = 0.097427
$ws.Cells.Item(13, 9).Value = 0.004586631199016116
$ws.Cells.Item(13, 10).Value = 0.004586631199016115
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.313829
$ws.Cells.Item(13, 14).Value = 0.941487
$ws.Cells.Item(13, 15).Value = 0.04225477273462707
$ws.Cells.Item(13, 16).Value = 0.04225477273462708
$ws.Cells.Item(13, 17).Value = 0.01019180599433333
$ws.Cells.Item(13, 18).Value = 0.091726253949
$ws.Cells.Item(13, 19).Value = 0.000193807058931976
$ws.Cells.Item(13, 20).Value = 0.000193807058931976
